$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter the Monday hours for the week commencing 43157 (row 7)
$ws.Range("B7").Value = 4

# Update the active selection to match the author's last edit position
$ws.Range("B7").Select()
